# Add a new "2022-Q3" sheet right after "总计" and before the existing
# "2022-Q2" sheet, populate it with the Q3 fund-holding data, and insert
# a matching summary row into the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for "2022-Q3" and bump the index
#    column for the rows that shift down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Copy formatting (style/border/font) from the row that is now row 3
# (the old row 2) onto the freshly inserted row 2, restricted to the
# A:D block actually used by this sheet so we don't smear a style
# across every column in the row.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 24
$total.Cells.Item(2, 4).Value = 1.67

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $q3.Cells.Item(1, $col + 2)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @('501079','大成科创主题混合（LOF）A','10.00','85.01','4.27','0.4270',8),
    @('012473','大成成长回报六个月持有混合A','7.48','80.86','4.01','0.2999',8),
    @('160921','大成多策略混合（LOF）A','2.90','87.66','7.52','0.2181',2),
    @('200006','长城消费增值混合','5.46','90.90','3.14','0.1714',10),
    @('002938','中银证券健康产业灵活配置混合','1.98','92.72','6.28','0.1243',2),
    @('010371','大成成长进取混合A','3.61','80.75','3.26','0.1177',10),
    @('016062','大成多策略混合（LOF）C','1.37','87.66','7.52','0.1030',2),
    @('014271','大成北交所两年定开混合A','3.45','65.31','1.70','0.0586',10),
    @('010372','大成成长进取混合C','1.52','80.75','3.26','0.0496',10),
    @('010434','红土创新医疗保健股票','0.44','94.49','5.02','0.0221',6),
    @('005044','国寿安保健康科学混合C','0.75','87.70','2.67','0.0200',8),
    @('012474','大成成长回报六个月持有混合C','0.37','80.86','4.01','0.0148',8),
    @('014272','大成北交所两年定开混合C','0.82','65.31','1.70','0.0139',10),
    @('519969','长信新利灵活配置混合','0.49','89.30','2.42','0.0119',8),
    @('001318','东方新策略灵活配置混合A','0.39','36.24','2.05','0.0080',1),
    @('400020','东方成长回报平衡混合','0.15','44.26','2.43','0.0036',1),
    @('005043','国寿安保健康科学混合A','0.08','87.70','2.67','0.0021',8),
    @('001657','长安鑫富领先灵活配置混合','0.06','49.60','2.51','0.0015',9),
    @('002060','东方新策略灵活配置混合C','0.06','36.24','2.05','0.0012',1),
    @('016198','大成科创主题混合（LOF）C','0.01','85.01','4.27','0.0004',8),
    @('001495','东方新价值混合A','0.01','29.03','2.22','0.0002',2),
    @('004166','东方价值挖掘灵活配置混合A','0.00','29.72','1.93',0,2),
    @('007686','东方价值挖掘灵活配置混合C','0.00','29.72','1.93',0,2),
    @('002162','东方新价值混合C','0.00','29.03','2.22',0,2)
)

$rowIndex = 2
foreach ($r in $rows) {
    $q3.Cells.Item($rowIndex, 1).Value = $rowIndex - 2

    # Fund code / name / scale / position / ratio are kept as literal
    # text (mirrors the source data, which stores these as strings,
    # e.g. "10.00" rather than the number 10).
    $q3.Cells.Item($rowIndex, 2).Value = "'" + $r[0]
    $q3.Cells.Item($rowIndex, 3).Value = "'" + $r[1]
    $q3.Cells.Item($rowIndex, 4).Value = "'" + $r[2]
    $q3.Cells.Item($rowIndex, 5).Value = "'" + $r[3]
    $q3.Cells.Item($rowIndex, 6).Value = "'" + $r[4]

    if ($r[5] -is [string]) {
        $q3.Cells.Item($rowIndex, 7).Value = "'" + $r[5]
    } else {
        $q3.Cells.Item($rowIndex, 7).Value = $r[5]
    }

    $q3.Cells.Item($rowIndex, 8).Value = $r[6]

    $rowIndex++
}
